$d = $word.ActiveDocument

# 1. Remove spurious "would" before "has preregistered"
$d.Content.Find.Execute("The user of the app would has preregistered", $true, $false, $false, $false, $false, $true, 1, $false, "The user of the app has preregistered", 2)

# 2. "with the red status" -> "to the red status"
$d.Content.Find.Execute("maps the user ID with the red status", $true, $false, $false, $false, $false, $true, 1, $false, "maps the user ID to the red status", 2)

# 3. "send packets to all users who had / have proximity" -> "sends packets to all users who have had proximity"
$d.Content.Find.Execute("send packets to all users who had / have proximity", $true, $false, $false, $false, $false, $true, 1, $false, "sends packets to all users who have had proximity", 2)

# 4. "this notification is the proximity" -> "this notification are proximity"
$d.Content.Find.Execute("this notification is the proximity", $true, $false, $false, $false, $false, $true, 1, $false, "this notification are proximity", 2)

# 5. "This notifies them they have been in contact and might be infected" -> "This notifies them that they have been in contact and may be infected"
$d.Content.Find.Execute("This notifies them they have been in contact and might be infected", $true, $false, $false, $false, $false, $true, 1, $false, "This notifies them that they have been in contact and may be infected", 2)

# 6. Drop leading space and "the " before COVID-19
$d.Content.Find.Execute(" Once the dashboard has been updated, all users in the area of the infection are alerted that someone near them has contracted the COVID-19.", $true, $false, $false, $false, $false, $true, 1, $false, "Once the dashboard has been updated, all users in the area of the infection are alerted that someone near them has contracted COVID-19.", 2)

# 7. "the database acknowledge that" -> "the database will acknowledge that" while reproducing the
#    run-split / _GoBack bookmark Word leaves behind after an in-place insertion edit.
$r = $d.Content
$r.Find.Execute("database acknowledge", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dbRange = $d.Range($r.Start, $r.Start + 8)
$dbRange.InsertAfter(" will")
$willRange = $d.Range($r.Start + 8, $r.Start + 13)
$bm = $d.Bookmarks.Add("_GoBack", $willRange)
$bmRange = $bm.Range
$bmRange.MoveStart(1, 5)
$d.Bookmarks.Add("_GoBack", $bmRange)
